# Add "That Will Never Work" by Marc Randolph to the finished book list
# (Sheet1 row 4): Title, Author, Start Date, Finish Date, Tags.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / Author / Tags are plain text values.
$ws.Range("A4").Value = "That Will Never Work"
$ws.Range("B4").Value = "Marc Randolph"
$ws.Range("E4").Value = "netflix;startups;business;ipo;technology"

# Start Date (1/3/2020) and Finish Date (1/5/2020) need to keep the same
# date number-format style already used by the other rows (style index
# reused, not a freshly minted number format) - copy the format down from
# the row above, then set the underlying date serial values.
$ws.Range("C3").Copy($ws.Range("C4")) | Out-Null
$ws.Range("C4").Value = 43833

$ws.Range("D3").Copy($ws.Range("D4")) | Out-Null
$ws.Range("D4").Value = 43835

# Move the active selection down to the next empty row, like Excel does
# after typing a new row of data.
$ws.Range("A5").Select() | Out-Null
